# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados" timestamp string in A1 (08:16 -> 08:46)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 08:46"

# Swap the Huelva / Huesca rows (row 53 <-> row 54).
# Row 53 currently holds Huelva's data, row 54 holds Huesca's data.
# After the edit, row 53 should hold Huesca's data and row 54 Huelva's data.
$newRow53 = New-Object 'object[,]' 1,5
$newRow53[0,0] = "Huesca"
$newRow53[0,1] = 37
$newRow53[0,2] = 0
$newRow53[0,3] = 37
$newRow53[0,4] = 0

$newRow54 = New-Object 'object[,]' 1,5
$newRow54[0,0] = "Huelva"
$newRow54[0,1] = 37
$newRow54[0,2] = 72
$newRow54[0,3] = 37
$newRow54[0,4] = 0

$ws.Range("A53:E53").Value = $newRow53
$ws.Range("A54:E54").Value = $newRow54
